$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(4493, 4558, 4782, 4872, 4872, 4929, 4935, 4979, 5068, 5132, 5132, 5215, 5215, 5225)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $newValues[$i]
}
